# Generate Report for handback
# Adds two new handback entries (5e269515-... and f424183b-...) to the
# Overview sheet as well as the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# New file identities being handed back
# ---------------------------------------------------------------------
$items = @(
    @{
        Id          = "5e269515-6661-4254-a937-700a2b62c7a6"
        Hash        = "72e5593b8f95246668dd5130427830932ed771b6"
        ZhHandoffDt = "2016-01-20 03:59:01"
        ZhHandbkDt  = "2016-01-20 03:59:48"
        DeHandoffDt = "2016-01-20 03:59:12"
        DeHandbkDt  = "2016-01-20 04:00:08"
    },
    @{
        Id          = "f424183b-efe1-469e-8484-6442cb9007d3"
        Hash        = "6923cd643c064e19feb5556a571141b1c741ad91"
        ZhHandoffDt = "2016-01-20 03:59:01"
        ZhHandbkDt  = "2016-01-20 03:59:48"
        DeHandoffDt = "2016-01-20 03:59:12"
        DeHandbkDt  = "2016-01-20 04:00:08"
    }
)

# ---------------------------------------------------------------------
# Overview sheet (columns: A=File Name, B=zh-cn, C=de-de)
# ---------------------------------------------------------------------
$overviewRow = 6
foreach ($item in $items) {
    $mdName = "$($item.Id).md"
    $mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName"

    $overview.Cells.Item($overviewRow, 1).Value = $mdName
    $overview.Cells.Item($overviewRow, 1).Style = "HyperLink"
    $overview.Hyperlinks.Add($overview.Cells.Item($overviewRow, 1), $mdUrl, "", "", $mdName)

    $overview.Cells.Item($overviewRow, 2).Value = $status
    $overview.Cells.Item($overviewRow, 3).Value = $status

    $overviewRow = $overviewRow + 1
}

# ---------------------------------------------------------------------
# Helper to populate a language detail sheet (zh-cn / de-de)
# Columns: A=Source File Name, B=Status, C=Correspond Handoff File,
#          D=Correspond Handoff Datetime, E=Target File,
#          F=Correspond Handback File, G=Correspond Handback DateTime,
#          H=Handoff Reason
# ---------------------------------------------------------------------
function Fill-LangSheet($ws, $lang, $handoffDtField, $handbkDtField) {
    $row = 6
    foreach ($item in $items) {
        $mdName  = "$($item.Id).md"
        $mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName"
        $xlfName = "$($item.Id).$($item.Hash).$lang.xlf"
        $xlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.$lang/xinjiang/$xlfName"
        $handoffDt = $item[$handoffDtField]
        $handbkDt  = $item[$handbkDtField]

        $ws.Cells.Item($row, 1).Value = $mdName
        $ws.Cells.Item($row, 1).Style = "HyperLink"
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $mdUrl, "", "", $mdName)

        $ws.Cells.Item($row, 2).Value = $status

        $ws.Cells.Item($row, 3).Value = $xlfName
        $ws.Cells.Item($row, 3).Style = "HyperLink"
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 3), $xlfUrl, "", "", $xlfName)

        $ws.Cells.Item($row, 4).Value = $handoffDt
        $ws.Cells.Item($row, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

        $ws.Cells.Item($row, 5).Value = $mdName
        $ws.Cells.Item($row, 5).Style = "HyperLink"
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $mdUrl, "", "", $mdName)

        $ws.Cells.Item($row, 6).Value = $xlfName
        $ws.Cells.Item($row, 6).Style = "HyperLink"
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $xlfUrl, "", "", $xlfName)

        $ws.Cells.Item($row, 7).Value = $handbkDt

        $ws.Cells.Item($row, 8).Value = "Include"

        $row = $row + 1
    }
}

Fill-LangSheet $zhcn "zh-cn" "ZhHandoffDt" "ZhHandbkDt"
Fill-LangSheet $dede "de-de" "DeHandoffDt" "DeHandbkDt"

Write-Host "Handback rows added."
